# Updated cryptos list with GitHub Actions
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = "97.129.97"
$ws.Range("E2").Value = "  +0.74%  "

# Row 3 (Ethereum)
$ws.Range("D3").Value = "3.691.97"
$ws.Range("E3").Value = "  +0.86%  "

# Row 4 (TetherUSD)
$ws.Range("E4").Value = "  -0.03%  "

# Row 5 (Solana)
$ws.Range("D5").Value = "'237.33"
$ws.Range("E5").Value = "  -1.79%  "

# Row 6 (XRP)
$ws.Range("E6").Value = "  +1.68%  "

# Row 7 (BNB)
$ws.Range("D7").Value = "'658.85"
$ws.Range("E7").Value = "  -0.04%  "

# Row 8 (Dogecoin)
$ws.Range("D8").Value = "'0.425"
$ws.Range("E8").Value = "  +0.53%  "

# Row 9 (Cardano)
$ws.Range("E9").Value = "  -1.41%  "

# Row 10 (USDC)
$ws.Range("E10").Value = "  -0.03%  "

# Row 11 (LidoStakedEther)
$ws.Range("D11").Value = "3.690.93"
$ws.Range("E11").Value = "  +0.86%  "

# Row 12 (was Avalanche -> TRON)
$ws.Range("B12").Value = "TRON"
$ws.Range("C12").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D12").Value = "'0.209"
$ws.Range("E12").Value = "  +2.53%  "

# Row 13 (was TRON -> ShibaInu)
$ws.Range("B13").Value = "ShibaInu"
$ws.Range("C13").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D13").Value = "'0.0000307"
$ws.Range("E13").Value = "  +13.31%  "

# Row 14 (was ShibaInu -> Avalanche)
$ws.Range("B14").Value = "Avalanche"
$ws.Range("C14").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D14").Value = "'44.13"
$ws.Range("E14").Value = "  -1.98%  "

# Row 15 (Toncoin)
$ws.Range("D15").Value = "'6.78"
$ws.Range("E15").Value = "  +0.74%  "

# Row 16 (WrappedliquidstakedEther2.0)
$ws.Range("D16").Value = "4.381.10"
$ws.Range("E16").Value = "  +0.90%  "

# Row 17 (WrappedBTC)
$ws.Range("D17").Value = "96.865.93"
$ws.Range("E17").Value = "  +0.72%  "

# Row 18 (Polkadot)
$ws.Range("D18").Value = "'9.19"
$ws.Range("E18").Value = "  +3.33%  "

# Row 19 (WrappedEther)
$ws.Range("D19").Value = "3.706.66"
$ws.Range("E19").Value = "  +2.27%  "

# Row 20 (Uniswap)
$ws.Range("D20").Value = "'13.01"
$ws.Range("E20").Value = "  +2.25%  "

# Row 21 (Chainlink)
$ws.Range("D21").Value = "'18.73"
$ws.Range("E21").Value = "  +2.66%  "

# Row 22 (Stellar)
$ws.Range("D22").Value = "'0.509"
$ws.Range("E22").Value = "  -3.88%  "

# Row 23 (BitcoinCash)
$ws.Range("D23").Value = "'520.76"
$ws.Range("E23").Value = "  -0.21%  "

# Row 24 (SuiNetwork)
$ws.Range("D24").Value = "'3.44"
$ws.Range("E24").Value = "  +0.25%  "

# Row 25 (PEPE)
$ws.Range("D25").Value = "'0.0000211"
$ws.Range("E25").Value = "  +3.34%  "

# Row 26 (NEARProtocol)
$ws.Range("D26").Value = "'6.95"
$ws.Range("E26").Value = "  +0.88%  "

# Row 27 (Hedera)
$ws.Range("E27").Value = "  +19.29%  "

# Row 28 (Litecoin)
$ws.Range("D28").Value = "'101.50"
$ws.Range("E28").Value = "  -0.39%  "

# Row 29 (Aptos)
$ws.Range("D29").Value = "'13.48"
$ws.Range("E29").Value = "  +3.75%  "

# Row 30 (InternetComputer(DFINITY))
$ws.Range("D30").Value = "'12.55"
$ws.Range("E30").Value = "  +1.55%  "

# Row 31 (PancakeSwap)
$ws.Range("D31").Value = "'3.03"
$ws.Range("E31").Value = "  -0.10%  "

# Row 33 (Cronos)
$ws.Range("D33").Value = "'0.191"
$ws.Range("E33").Value = "  +2.70%  "

# Row 34 (Fetch.AI)
$ws.Range("D34").Value = "'1.88"
$ws.Range("E34").Value = "  +1.85%  "

# Row 35 (Binance-PegBSC-USD)
$ws.Range("D35").Value = "'0.999"
$ws.Range("E35").Value = "  -1.66%  "

# Row 36 (EthereumClassic)
$ws.Range("E36").Value = "  -2.26%  "

# Row 37 (Bittensor)
$ws.Range("D37").Value = "'651.09"
$ws.Range("E37").Value = "  +3.16%  "

# Row 38 (PolygonEcosystemToken)
$ws.Range("D38").Value = "'0.593"
$ws.Range("E38").Value = "  +0.97%  "

# Row 39 (RenderToken)
$ws.Range("D39").Value = "'8.82"
$ws.Range("E39").Value = "  +1.05%  "

# Row 41 (was Filecoin -> Algorand)
$ws.Range("B41").Value = "Algorand"
$ws.Range("C41").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D41").Value = "'0.503"
$ws.Range("E41").Value = "  +12.44%  "

# Row 42 (was Algorand -> Filecoin)
$ws.Range("B42").Value = "Filecoin"
$ws.Range("C42").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D42").Value = "'6.82"
$ws.Range("E42").Value = "  +8.72%  "

# Row 43 (ImmutableX)
$ws.Range("E43").Value = "  +4.86%  "

# Row 44 (Kaspa)
$ws.Range("E44").Value = "  +1.26%  "

# Row 45 (EnergySwap)
$ws.Range("D45").Value = "'40.44"
$ws.Range("E45").Value = "  -9.80%  "

# Row 46 (ARBITRUM)
$ws.Range("D46").Value = "'0.962"
$ws.Range("E46").Value = "  +0.39%  "

# Row 47 (VeChain)
$ws.Range("D47").Value = "'0.0467"
$ws.Range("E47").Value = "  +1.50%  "

# Row 48 (Stacks)
$ws.Range("D48").Value = "'2.29"
$ws.Range("E48").Value = "  +0.19%  "

# Row 49 (WhiteBITCoin)
$ws.Range("D49").Value = "'23.63"
$ws.Range("E49").Value = "  +0.03%  "

# Row 50 (Cosmos)
$ws.Range("D50").Value = "'8.75"
$ws.Range("E50").Value = "  +2.33%  "

# Row 51 (MantraDAO)
$ws.Range("D51").Value = "'3.52"
$ws.Range("E51").Value = "  -0.45%  "
